$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 5, pushing the existing rows 5-17 down to 6-18 ---
$ws.Range("A5").EntireRow.Insert()

# --- Fix / normalize a few existing labels (spacing added) ---
$ws.Range("D4").Value  = "Socialización Gerencia"
$ws.Range("D6").Value  = "Implementación Last Planner"
$ws.Range("D7").Value  = "Implementación Grilla LP"
$ws.Range("D10").Value = "Seguimiento Y Control"
$ws.Range("D11").Value = "Desarrollo Software"

# --- Populate the newly inserted row 5 ---
$ws.Range("A5").Value = "1.1.1.1."
$ws.Range("B5").Value = "1.1.1."
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "Presentación socialización mejoras Gerencia"
$ws.Range("E5").Value = "Pptx"
$ws.Range("F5").Value = "Propuesta de mejoras"
$ws.Range("G5").Value = "https://grupomarval-my.sharepoint.com/personal/cpulgarin_marval_com_co/Documents/Datos%20adjuntos/PropuestaMejoraIngenieria.pptx?web=1"

# --- Turn the URL in G5 into a real hyperlink (Excel auto-creates the Hyperlink style) ---
$ws.Hyperlinks.Add($ws.Range("G5"), "https://grupomarval-my.sharepoint.com/personal/cpulgarin_marval_com_co/Documents/Datos%20adjuntos/PropuestaMejoraIngenieria.pptx?web=1")

# --- Grow the table so the new row (and the trailing blank row) are included ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:H18"))

# --- Re-fit the columns whose best-fit width needs to grow for the new content ---
$ws.Columns.Item(3).EntireColumn.AutoFit()
$ws.Columns.Item(4).EntireColumn.AutoFit()
$ws.Columns.Item(7).EntireColumn.AutoFit()

# --- Restore selection to the newly-edited cell ---
$ws.Range("D5").Select()
